# Weekly update: insert a new price observation as the first (most recent)
# record for this market/variety subset, pushing the existing historical
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 136 - this shifts the existing rows
# 136..164 down to 137..165 and grows the used range to A1:R165.
$ws.Rows.Item(136).Insert()

# Populate the new row 136 with this week's observation.
$ws.Cells.Item(136, 1).Value = 4
$ws.Cells.Item(136, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(136, 3).Value = "Los Lagos"
$ws.Cells.Item(136, 4).Value = 44543
$ws.Cells.Item(136, 5).Value = 10
$ws.Cells.Item(136, 6).Value = 100112032
$ws.Cells.Item(136, 7).Value = "Zapallo italiano"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 70
$ws.Cells.Item(136, 11).Value = 11000
$ws.Cells.Item(136, 12).Value = 11000
$ws.Cells.Item(136, 13).Value = 11000
$ws.Cells.Item(136, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(136, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(136, 16).Value = 220
$ws.Cells.Item(136, 17).Value = 50
$ws.Cells.Item(136, 18).Value = "Hortaliza"
